$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175 (pushes existing row 175 and everything
# below it down by one, growing the sheet from A1:R278 to A1:R279).
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the latest week's data.
$ws.Cells.Item(175, 1).Value = 10
$ws.Cells.Item(175, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(175, 3).Value = "La Araucanía"
$ws.Cells.Item(175, 4).Value = 45029
$ws.Cells.Item(175, 5).Value = 9
$ws.Cells.Item(175, 6).Value = 100114007
$ws.Cells.Item(175, 7).Value = "Jengibre"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 80
$ws.Cells.Item(175, 11).Value = 25000
$ws.Cells.Item(175, 12).Value = 25000
$ws.Cells.Item(175, 13).Value = 25000
$ws.Cells.Item(175, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(175, 15).Value = "Perú"
$ws.Cells.Item(175, 16).Value = 1923
$ws.Cells.Item(175, 17).Value = 13
$ws.Cells.Item(175, 18).Value = "Hortaliza"
